# Update Target Depth Data for Week 15 (logged) and Week 16 (simulated)
$wb = $excel.ActiveWorkbook

# OFF sheet (row 3: Road totals)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 408
$wsOff.Range("C3").Value = 298
$wsOff.Range("D3").Value = 110
$wsOff.Range("E3").Value = 54

# DEF sheet (row 3: Road totals)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 444
$wsDef.Range("C3").Value = 303
$wsDef.Range("D3").Value = 100
$wsDef.Range("E3").Value = 38
